# serviceapplication.xlsx - "done with the modification of the excelsheet and upload data"
#
# Row 2's "Marine_Hod_Approval" column is being repurposed/split: the header
# (N1) becomes "Hod_Approval" while a new "inspection type" column (C) gets
# populated with "mechanical" (row 2) / "manual" (row 3) values, row 2's
# inspection date moves back a month, and a brand-new row 3 is appended with
# mostly-1 data (fee-verified column F intentionally left blank, matching the
# source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 edits -----------------------------------------------------------
# Service Type ID becomes a text classification instead of a bare number.
$ws.Range("C2").Value = "mechanical"
# Date_Of_Inspection moves from 2023-02-01 (44958) to 2023-01-01 (44928).
$ws.Range("E2").Value = 44928

# --- New row 3 --------------------------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "manual"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 44928
# Keep the same custom date format as E2 so it shares the existing style
# instead of minting a near-duplicate number format.
$ws.Range("E3").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
# F3 (fee-verified) is intentionally left empty for this row.
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1

# --- Header rename ----------------------------------------------------------
# Do this last so the new shared-string table ends up ordered the same way
# as the target workbook: mechanical, manual, Hod_Approval.
$ws.Range("N1").Value = "Hod_Approval"

# --- View/selection state ----------------------------------------------------
$ws.Range("O3").Select() | Out-Null
